$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet2 (Feuil2): add the new "distance -> sqrt" helper columns G:M for
# rows 16-22, matching the POWER() grid already present in V:AB.
# ---------------------------------------------------------------------------

# Row 16: G16 is a single formula, H16:M16 is one shared-formula group.
$ws2.Range("G16").Formula = "=SQRT(V16)"
$ws2.Range("H16:M16").Formula = "=SQRT(W16)"

# Rows 17-22: each column G-M becomes its own shared-formula group running
# down the rows (mirrors how V17:V22, W17:W22, etc. were originally built).
$ws2.Range("G17:G22").Formula = "=SQRT(V17)"
$ws2.Range("H17:H22").Formula = "=SQRT(W17)"
$ws2.Range("I17:I22").Formula = "=SQRT(X17)"
$ws2.Range("J17:J22").Formula = "=SQRT(Y17)"
$ws2.Range("K17:K22").Formula = "=SQRT(Z17)"
$ws2.Range("L17:L22").Formula = "=SQRT(AA17)"
$ws2.Range("M17:M22").Formula = "=SQRT(AB17)"

# New column T (20) gets an explicit (custom) width, same value as the
# sheet's default column width.
$ws2.Columns.Item(20).ColumnWidth = 3.2

# ---------------------------------------------------------------------------
# Conditional formatting: the old color-scale on V16:AB22 is replaced with a
# "cell value < 12.26" rule, and a new "cell value < 3.501" rule is added on
# the new G16:M22 range. Both rules use the same red-text / pink-fill dxf.
#
# The new G16:M22 rule ends up listed before V16:AB22 in the saved XML (so it
# is added first), but V16:AB22 keeps the lower priority number (1, vs 2 for
# G16:M22) and the lower dxf index (0, vs 1) - so the dxf style is applied to
# V16:AB22 first even though its conditionalFormatting block is added second.
# ---------------------------------------------------------------------------
$rngNew = $ws2.Range("G16:M22")
$fcNew = $rngNew.FormatConditions.Add(1, 6, "3.501")

$rngOld = $ws2.Range("V16:AB22")
$rngOld.FormatConditions.Delete()
$fcOld = $rngOld.FormatConditions.Add(1, 6, "12.26")

$fcOld.Font.Color = 393372
$fcOld.Interior.Color = 13551615
$fcNew.Font.Color = 393372
$fcNew.Interior.Color = 13551615

$fcOld.Priority = 1
$fcNew.Priority = 2

# ---------------------------------------------------------------------------
# Sheet view / selection changes.
#   Feuil2 loses tabSelected and its selection moves from V16:AB22 to T16.
#   Feuil1 gains tabSelected (becomes the active sheet) with selection D13,
#   scrolled so row 12 is at the top.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("T16").Select()

$ws1.Activate()
$ws1.Range("D13").Select()
$excel.ActiveWindow.ScrollRow = 12
